$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = "extinct (post 1500)"
$ws.Range("B12").Value = "absent"

$ws.Range("A13").Value = "absent"
$ws.Range("B13").Value = "absent"

$ws.Range("A14").Value = "present"
$ws.Range("B14").Value = "present"

$ws.Range("A15").Value = "unknown"
$ws.Range("B15").Value = "uncertain"

$ws.Range("A16").Select()
